$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I3:I14").Select()
$ws.Range("I4").Select($false)
Write-Output "done"
